# Auto-generated edit script for 杭州-漫展信息.xlsx update
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (exhibitions): update "want to go" counts (F column) ---
$ws1.Range("F3").Value = 2500
$ws1.Range("F5").Value = 905
$ws1.Range("F7").Value = 1312
$ws1.Range("F8").Value = 1673
$ws1.Range("F9").Value = 173
$ws1.Range("F11").Value = 2327
$ws1.Range("F18").Value = 8483
$ws1.Range("F20").Value = 6542
$ws1.Range("F21").Value = 10575
$ws1.Range("F23").Value = 184
$ws1.Range("F24").Value = 202
$ws1.Range("F29").Value = 82
$ws1.Range("F30").Value = 24
$ws1.Range("F31").Value = 12
$ws1.Range("F32").Value = 4435
$ws1.Range("F33").Value = 334
$ws1.Range("F34").Value = 422

# --- Sheet "演出" (performances): update "want to go" counts (F column) ---
$ws2.Range("F19").Value = 4
$ws2.Range("F20").Value = 13
$ws2.Range("F23").Value = 4

# --- Sheet "全部类型" (all types): update "want to go" counts (F column) ---
$ws4.Range("F6").Value = 2500
$ws4.Range("F8").Value = 905
$ws4.Range("F10").Value = 1312
$ws4.Range("F12").Value = 1673
$ws4.Range("F14").Value = 173
$ws4.Range("F15").Value = 2327
$ws4.Range("F24").Value = 8483
$ws4.Range("F26").Value = 6542
$ws4.Range("F27").Value = 10575
$ws4.Range("F30").Value = 184
$ws4.Range("F31").Value = 202
$ws4.Range("F40").Value = 4435
$ws4.Range("F46").Value = 4
$ws4.Range("F47").Value = 422

# --- Sheet "全部类型": rows 33-39 content shifts (one old event removed, one new event inserted) ---
$ws4.Range("C33").Value = '杭州·ET金色齿轮国乙同人only'
$ws4.Range("D33").Value = '转塘街道珊瑚沙东路9号 杭州白金汉爵大酒店'
$ws4.Range("E33").Value = '2024.11.16 09:30-11.16 22:00'
$ws4.Range("F33").Value = 519
$ws4.Range("G33").Value = 25
$ws4.Range("H33").Value = 'https://show.bilibili.com/platform/detail.html?id=92511'
$ws4.Range("I33").Value = '//i1.hdslb.com/bfs/openplatform/202409/XfT00A611726134427042.jpeg'

$ws4.Range("C34").Value = '杭州·“法国姐姐”乔伊丝·乔纳森《小意思》2024巡回演唱会【特邀嘉宾陈丽君】'
$ws4.Range("D34").Value = '杭州市西湖区省府路9号 浙江省人民大会堂'
$ws4.Range("E34").Value = '2024.11.16 19:30-11.16 21:30'
$ws4.Range("F34").Value = 5
$ws4.Range("G34").Value = 280
$ws4.Range("H34").Value = 'https://show.bilibili.com/platform/detail.html?id=92078'
$ws4.Range("I34").Value = '//i2.hdslb.com/bfs/openplatform/202409/AE6VYTdf1725614295764.jpeg'

$ws4.Range("B35").Value = '2024-11-22'
$ws4.Range("C35").Value = '杭州·【早鸟5折起】《LALALAND爱乐之城》浪漫主题音乐会'
$ws4.Range("D35").Value = '曙光路31号 浙江音乐厅'
$ws4.Range("E35").Value = '2024.11.22 19:30-11.22 21:00'
$ws4.Range("F35").Value = 8
$ws4.Range("G35").Value = 100
$ws4.Range("H35").Value = 'https://show.bilibili.com/platform/detail.html?id=92725'
$ws4.Range("I35").Value = '//i2.hdslb.com/bfs/openplatform/202409/97ZKC3qG1727059280000.jpeg'

$ws4.Range("B36").Value = '2024-11-23'
$ws4.Range("C36").Value = '杭州·奥斯卡·罗曼耶卓（O叔）钢琴独奏音乐会'
$ws4.Range("D36").Value = '建国南路280号 杭州红星剧院'
$ws4.Range("E36").Value = '2024.11.23 19:30-11.23 21:15'
$ws4.Range("F36").Value = 12
$ws4.Range("G36").Value = 180
$ws4.Range("H36").Value = 'https://show.bilibili.com/platform/detail.html?id=91019'
$ws4.Range("I36").Value = '//i2.hdslb.com/bfs/openplatform/202408/PJRlAXdp1724126404150.jpeg'

$ws4.Range("C37").Value = '杭州·火影忍者同人only2.0 日夜连场'
$ws4.Range("D37").Value = '金一路79号 XPACE湾区数字公园'
$ws4.Range("E37").Value = '2024.11.23 10:00-11.23 22:30'
$ws4.Range("F37").Value = 188
$ws4.Range("G37").Value = 69
$ws4.Range("H37").Value = 'https://show.bilibili.com/platform/detail.html?id=92097'
$ws4.Range("I37").Value = '//i2.hdslb.com/bfs/openplatform/202409/q3I7lKmY1725591212982.jpeg'

$ws4.Range("C38").Value = '杭州·相聚广陵代号鸢同人only3.0-三千世界'
$ws4.Range("D38").Value = '康候圣街99号 顺丰创新中心'
$ws4.Range("E38").Value = '2024.11.23 09:30-11.23 17:00'
$ws4.Range("F38").Value = 165
$ws4.Range("G38").Value = 80
$ws4.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=92672'
$ws4.Range("I38").Value = '//i0.hdslb.com/bfs/openplatform/202409/NR40ECNZ1726740199589.jpeg'

$ws4.Range("B39").Value = '2024-12-07'
$ws4.Range("C39").Value = '杭州·eva同人only'
$ws4.Range("D39").Value = '康候圣街99号 顺丰创新中心'
$ws4.Range("E39").Value = '2024.12.07 10:00-12.07 17:30'
$ws4.Range("F39").Value = 12
$ws4.Range("G39").Value = 79
$ws4.Range("H39").Value = 'https://show.bilibili.com/platform/detail.html?id=93190'
$ws4.Range("I39").Value = '//i1.hdslb.com/bfs/openplatform/202409/onI1p61S1727690051879.png'

# --- Sheet "全部类型": rows 48-49 content shift (old row48 event dropped, new event added at row49) ---
$ws4.Range("B48").Value = '2025-01-01'
$ws4.Range("C48").Value = '杭州·【早鸟优惠】大型正版授权互动卡通儿童剧《海底小纵队之深海探秘》'
$ws4.Range("D48").Value = '湖墅南路136-138号 浙话艺术剧院'
$ws4.Range("E48").Value = '2025.01.01 10:30-01.01 11:40'
$ws4.Range("F48").Value = 1
$ws4.Range("G48").Value = 40
$ws4.Range("H48").Value = 'https://show.bilibili.com/platform/detail.html?id=92951'
$ws4.Range("I48").Value = '//i2.hdslb.com/bfs/openplatform/202409/oZlaKX931727335820196.jpeg'

$ws4.Range("B49").Value = '2025-01-08'
$ws4.Range("C49").Value = '杭州·【早鸟限时六五折】维也纳皇家交响乐团2025新年音乐会'
$ws4.Range("D49").Value = '桥弄街399号（运河中央公园附近） 杭州运河大剧院'
$ws4.Range("E49").Value = '2025.01.08 19:30-01.08 21:00'
$ws4.Range("F49").Value = 4
$ws4.Range("G49").Value = 312
$ws4.Range("H49").Value = 'https://show.bilibili.com/platform/detail.html?id=92877'
$ws4.Range("I49").Value = '//i1.hdslb.com/bfs/openplatform/202409/dKSKfgEx1727240509662.png'

